# Applies the SEO.docx edit described in the commit:
# "Pongo contenido a las etiquetas meta de keywords."
#
# Summary of structural changes:
#  1. Split the "Reduzco tamaño de los elementos audiovisuales:" paragraph
#     so the "-video1.mp4" line becomes its own (indented) paragraph.
#  2. Indent the "-contacto.jpg" / "-index.jpg" paragraphs.
#  3. Trim the "Cambio de formato png a svg ... img (y corrijo los html ...)"
#     paragraph down to "... img:" and drop the _GoBack bookmark from it.
#  4. Append six new paragraphs (the continuation of the log), moving the
#     _GoBack bookmark to the very end of the document.
#  5. Indent the six "-logo-instagram" / "-Mail" / ... list paragraphs
#     (done last so the newly appended paragraphs above don't inherit it).

$d = $word.ActiveDocument

# --- 1. Split "Reduzco tamaño ... audiovisuales:" / "-video1.mp4" ---------
$r = $d.Paragraphs.Item(2).Range
$r.Find.Execute("audiovisuales:")
$r.Collapse(0)
$r.InsertParagraphAfter()

# New paragraph (now index 3) gets a left indent of 708 twips (= 35.4 pt).
$d.Paragraphs.Item(3).Range.ParagraphFormat.LeftIndent = 35.4

# --- 2. Indent "-contacto.jpg" and "-index.jpg" ---------------------------
$d.Paragraphs.Item(4).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(5).Range.ParagraphFormat.FirstLineIndent = 35.4

# --- 3. Trim the "Cambio de formato ..." paragraph and drop the bookmark --
$p8 = $d.Paragraphs.Item(8).Range
$p8.Find.Execute(" (y corrijo los html para que no se rompan los links a las imágenes referidas):", $true, $false, $false, $false, $false, $true, 1, $false, ":", 2)

if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}

# --- 4. Append the six new log paragraphs at the end of the document ------
# (Do this before re-indenting the list paragraphs below, so the new
# paragraphs do not inherit their first-line indent.)
$newParagraphs = @(
    "Corrijo los html para que no se rompan los links a las imágenes referidas.",
    "17/11/2021",
    "Pongo contenido a las etiquetas meta de keywords.",
    "Pongo contenido a las etiquetas meta de description.",
    "Corroboro que mis title son correctos y adecuados para un buen SEO.",
    "Corroboro que mi proyecto no cuenta con niveles de carpeta innecesarias."
)

$end = $d.Content
$end.Collapse(0)
for ($i = 0; $i -lt $newParagraphs.Length; $i++) {
    $end.InsertParagraphAfter()
    $end.Collapse(0)
    $newPara = $d.Paragraphs.Last
    $text = $newParagraphs[$i]
    if ($i -eq $newParagraphs.Length - 1) {
        # Append a throw-away placeholder character on the very last
        # paragraph so the bookmark-insertion point (below) never lands on
        # the document's absolute last position, then trim it back off.
        $text = $text + "Z"
    }
    $newPara.Range.InsertBefore($text)
    $end = $d.Content
    $end.Collapse(0)
}

$placeholderEnd = $d.Content.End
$bookmarkPos = $placeholderEnd - 2
$bookmarkRange = $d.Range($bookmarkPos, $bookmarkPos)
if ($d.Bookmarks.Exists("_GoBack")) {
    $d.Bookmarks.Item("_GoBack").Delete()
}
$d.Bookmarks.Add("_GoBack", $bookmarkRange)

$trimEnd = $d.Content.End
$placeholderRange = $d.Range($trimEnd - 2, $trimEnd - 1)
$placeholderRange.Delete()

# --- 5. Indent the six list paragraphs that precede the appended text -----
$d.Paragraphs.Item(9).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(10).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(11).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(12).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(13).Range.ParagraphFormat.FirstLineIndent = 35.4
$d.Paragraphs.Item(14).Range.ParagraphFormat.FirstLineIndent = 35.4
